$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.698.13'
$ws.Range("E2").Value = '  +1.53%  '

$ws.Range("D3").Value = '2.311.41'
$ws.Range("E3").Value = '  +0.34%  '

$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '319.39'
$ws.Range("E5").Value = '  +2.63%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.85'
$ws.Range("E6").Value = '  -1.14%  '

$ws.Range("E7").Value = '  +1.20%  '

$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.613'
$ws.Range("E9").Value = '  +1.41%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.87'
$ws.Range("E10").Value = '  +0.28%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0912'
$ws.Range("E11").Value = '  +0.32%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.38'
$ws.Range("E12").Value = '  +1.22%  '

$ws.Range("E13").Value = '  +1.03%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.971'
$ws.Range("E14").Value = '  +0.55%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.39'
$ws.Range("E15").Value = '  -0.03%  '

$ws.Range("D16").Value = '2.657.00'
$ws.Range("E16").Value = '  +0.06%  '

$ws.Range("D17").Value = '2.324.95'
$ws.Range("E17").Value = '  +0.27%  '

$ws.Range("D18").Value = '42.566.64'
$ws.Range("E18").Value = '  +1.11%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.54'
$ws.Range("E19").Value = '  -0.32%  '

$ws.Range("E20").Value = '  +1.57%  '

$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '283.04'
$ws.Range("E21").Value = '  +9.59%  '

$ws.Range("B22").Value = 'Litecoin'
$ws.Range("C22").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.33'
$ws.Range("E22").Value = '  -1.26%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.60'
$ws.Range("E23").Value = '  +3.96%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.95'
$ws.Range("E24").Value = '  +19.29%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.29'
$ws.Range("E25").Value = '  +1.11%  '

$ws.Range("E26").Value = '  -0.32%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.92'
$ws.Range("E27").Value = '  -0.31%  '

$ws.Range("E28").Value = '  +4.44%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '23.00'
$ws.Range("E29").Value = '  +1.25%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.21'
$ws.Range("E30").Value = '  +1.91%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '164.77'
$ws.Range("E31").Value = '  +0.37%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0879'
$ws.Range("E32").Value = '  -0.52%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.94'
$ws.Range("E33").Value = '  +1.99%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.138'
$ws.Range("E34").Value = '  +7.38%  '

$ws.Range("E35").Value = '  -9.43%  '

$ws.Range("E36").Value = '  +0.54%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0371'
$ws.Range("E37").Value = '  +6.18%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.63'
$ws.Range("E38").Value = '  +3.13%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.74'
$ws.Range("E39").Value = '  +2.65%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.80'
$ws.Range("E40").Value = '  +3.80%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.51'
$ws.Range("E41").Value = '  +2.87%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '99.06'
$ws.Range("E42").Value = '  +1.25%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '70.00'
$ws.Range("E43").Value = '  +0.49%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.228'
$ws.Range("E44").Value = '  -0.58%  '

$ws.Range("E45").Value = '  -0.07%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.21'
$ws.Range("E46").Value = '  +0.90%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '113.11'
$ws.Range("E47").Value = '  +2.02%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '79.55'
$ws.Range("E48").Value = '  +9.05%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.94'
$ws.Range("E49").Value = '  -0.37%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.31'
$ws.Range("E50").Value = '  -0.94%  '

$ws.Range("D51").Value = '1.609.30'
$ws.Range("E51").Value = '  +5.12%  '
